# Generate Report for Handoff
#
# "b.md" has been handed off again for zh-cn and de-de: update the Status
# to "Ready for handoff" on the Overview sheet and on each locale sheet,
# and record the new handoff file name / timestamp on each locale sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b.md (row 3) -> Status columns for zh-cn (B) and de-de (C)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for b.md (row 3)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-10 05:50:41"

# --- de-de sheet: row for b.md (row 3)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-10 05:50:51"
